$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 183 (before current row 184).
# This shifts existing rows 184:305 down to 186:307, matching the
# dimension growth from A1:R305 to A1:R307 seen in the diff.
$ws.Rows.Item(184).Resize(2).Insert()

# Populate the two newly inserted rows (184 and 185) with their data.
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg,
# F Categoría ID, G Categoría, H Variedad, I Calidad, J Volumen,
# K Precio mínimo, L Precio máximo, M Precio promedio ponderado,
# N Unidad de comercialización, O Origen, P Precio $/Kg,
# Q Kg o Unidades, R Clasificación

$row184 = @(10, "Vega Modelo de Temuco", "La Araucanía", 45176, 9, 100112012, "Espinaca", "Sin especificar", "Primera", 50, 10000, 10000, 10000, "`$/cuna 10 kilos", "Región Metropolitana", 1000, 10, "Hortaliza")
for ($i = 0; $i -lt $row184.Length; $i++) {
    $ws.Cells.Item(184, $i + 1).Value = $row184[$i]
}
$ws.Range("D184").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row185 = @(10, "Vega Modelo de Temuco", "La Araucanía", 45176, 9, 100112012, "Espinaca", "Sin especificar", "Primera", 60, 12000, 12000, 12000, "`$/docena de atados", "Región de La Araucanía", 4000, 3, "Hortaliza")
for ($i = 0; $i -lt $row185.Length; $i++) {
    $ws.Cells.Item(185, $i + 1).Value = $row185[$i]
}
$ws.Range("D185").NumberFormat = "YYYY-MM-DD HH:MM:SS"
